# Applies the ArmLab_P2 edit described by the commit:
#   - Adds a sentence about storing path points as joint-angle matrices to
#     the end of the "2.1" paragraph.
#   - Rewrites the "2.2" (Wavefront Planning) paragraph: merges the split
#     "T"/"he robot..." runs, wraps each "Wavefront" occurrence in
#     spell-check proofErr markers, replaces the "Besides, it was useful..."
#     tail with "We had also used Wavefront planning...", and removes the
#     trailing "the strategy to fit the arm-lab." run (folded into the
#     rewritten text above).
#   - Splits what used to be the start of "2.3" into its own paragraph with
#     expanded PID-control text, and moves the _GoBack bookmark to the end
#     of that new paragraph.

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, [string]$prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# --- Change 1: append a new run to the end of the "2.1" paragraph ---
$idx21 = Find-ParagraphIndex $d "2.1 The viable points"
$p21 = $d.Paragraphs($idx21)
$p21End = $p21.Range.End
$insertRange1 = $d.Range($p21End - 1, $p21End - 1)
$insertRange1.InsertAfter(" We will store the points along the path in a matrix as angles for each joint. This makes it easy for us to write our path following code since all the angles are easily accessible. ")

# --- Change 2: rewrite the "2.2" paragraph and split a new "2.3" paragraph out of it ---
$idx22 = Find-ParagraphIndex $d "2.2"
$idx23 = Find-ParagraphIndex $d "2.3"
$startPos = $d.Paragraphs($idx22).Range.Start
$endPos = $d.Paragraphs($idx23).Range.End
$fullRange = $d.Range($startPos, $endPos)

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>2.2</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">The robot will reach the goals using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Wavefront</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Planning. The biggest benefit to using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Wavefront</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for this lab is because the configuration created for part 1 made it easy to define the world in part2 by converting the configuration space to a binary grid. We defined ones for obstacles and zeros for reachable space, then performed the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Wavefront</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Planner. </w:t></w:r><w:r><w:t xml:space="preserve">We had also used </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Wavefront</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> planning for the Motion Planning lab, which made it easy for us to reuse a lot of our old code and adapt it to fit the arm lab. </w:t></w:r></w:p><w:p><w:r><w:t>2.3 Using two PID controls with feed forward; one PID con</w:t></w:r><w:r><w:t>trol</w:t></w:r><w:r><w:t xml:space="preserve"> for each link.</w:t></w:r><w:r><w:t xml:space="preserve"> We can calculate the error between our current base angle and the target base angle and use that to change our torque for the base joint. The same can be done for the second joint. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$fullRange.InsertXML($xmlFrag)

Write-Output $d.Content.Text
